$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Principal")
$ws.Activate()

$ws.Range("D2").Value = "OK"
$ws.Range("D3").Value = "OK"
$ws.Range("D4").Value = "OK"
$ws.Range("D5").Value = "OK"

$ws.Range("D4").Select()
